$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("K18").Value = 36.02
$ws1.Range("L25").Value = 855.36
$ws1.Range("Q35").Value = 834.48
$ws1.Range("K55").Value = "3 de 53"
$ws1.Range("L55").Value = "8 de 53"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F18").Value = 5548.51
$ws2.Range("F25").Value = 3201.65
$ws2.Range("F35").Value = 3622.99
$ws2.Range("F55").Value = 83050.90000000001

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 450.99
$ws3.Range("E12").Value = -100.99
$ws3.Range("F12").Value = 1.288542857142857

$ws3.Range("D14").Value = 1504.94
$ws3.Range("E14").Value = -538.9400000000001
$ws3.Range("F14").Value = 1.557908902691511

$ws3.Range("D15").Value = 22904.13
$ws3.Range("E15").Value = -9404.130000000001
$ws3.Range("F15").Value = 1.696602222222222

$ws3.Range("D19").Value = 85465.83
$ws3.Range("E19").Value = 8981.610645179149
$ws3.Range("F19").Value = 0.9049036100520571
